# Add three new diagnosis rows (14, 15, 16) used to check the
# "Low Bp" symptom handling, as described in the commit message
# "add lowbp symtoms to be checked".
#
# We seed each new row by copying an existing row that already has an
# empty "Underlying Health Issues" cell (row 5), so the empty cell is
# preserved in the sheet, and then overwrite every column with the
# actual values for the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 14 - Id 13 ----
$ws.Range("A5:Q5").Copy($ws.Range("A14:Q14"))
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Otra"
$ws.Range("C14").Value = "Baker"
$ws.Range("D14").Value = "robinsondave_876@yahoo.com"
$ws.Range("E14").Value = 113
$ws.Range("F14").Value = 34
$ws.Range("G14").Value = "Aches"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = $false
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = "Low Risk"

# ---- Row 15 - Id 14 (note the mistyped email, matching the source data) ----
$ws.Range("A5:Q5").Copy($ws.Range("A15:Q15"))
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Otra"
$ws.Range("C15").Value = "Baker"
$ws.Range("D15").Value = "robinsondave_876@yaho..com"
$ws.Range("E15").Value = 113
$ws.Range("F15").Value = 34
$ws.Range("G15").Value = "Aches"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = $false
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = "Low Risk"

# ---- Row 16 - Id 15 ----
$ws.Range("A5:Q5").Copy($ws.Range("A16:Q16"))
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Otra"
$ws.Range("C16").Value = "Baker"
$ws.Range("D16").Value = "robinsondave_876@yahoo.com"
$ws.Range("E16").Value = 113
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = "Aches"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = $false
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = "Low Risk"
